$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Version bump
$ws1.Range("B3").Value = "6.0.0"

# Publication date update
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value
$ws1.Range("B9").Value = "Alvearie Team"

# Remove the duplicate "Contact / No display for ContactDetail" row (row 11),
# which shifts everything below it up by one row.
$ws1.Rows.Item(11).Delete()

# Replace the remaining Contact row with Jurisdiction info
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# --- Sheet 2: "Elements" ---
$ws2 = $wb.Worksheets.Item("Elements")

# Root extension row: Short / Definition now reflect the specific extension
$ws2.Range("K2").Value = "Employee Termination Date"
$ws2.Range("L2").Value = "Last date of employment for the employee"
